$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "25.960.82"
$cell.ClearFormats()

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.744.80"
$cell.ClearFormats()

$ws.Range("E3").Value = "  -0.25%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.ClearFormats()

$ws.Range("E4").Value = "  +0.08%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "248.75"
$cell.ClearFormats()

$ws.Range("E5").Value = "  +5.00%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.ClearFormats()

$ws.Range("E6").Value = "  +0.05%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.5038"
$cell.ClearFormats()

$ws.Range("E7").Value = "  -4.87%  "

$ws.Range("E8").Value = "  -2.41%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06184"
$cell.ClearFormats()

$ws.Range("E9").Value = "  +0.06%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.07279"
$cell.ClearFormats()

$ws.Range("E10").Value = "  +1.41%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "1.744.93"
$cell.ClearFormats()

$ws.Range("E11").Value = "  -0.24%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.6529"
$cell.ClearFormats()

$ws.Range("E12").Value = "  +1.02%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "15.16"
$cell.ClearFormats()

$ws.Range("E13").Value = "  -1.99%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "4.647"
$cell.ClearFormats()

$ws.Range("E14").Value = "  +0.32%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "77.66"
$cell.ClearFormats()

$ws.Range("E15").Value = "  -1.00%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.ClearFormats()

$ws.Range("E16").Value = "  +0.06%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.9996"
$cell.ClearFormats()

$ws.Range("E17").Value = "  +0.08%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "25.983.68"
$cell.ClearFormats()

$ws.Range("E18").Value = "  -0.18%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "11.84"
$cell.ClearFormats()

$ws.Range("E19").Value = "  +0.40%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.000006827"
$cell.ClearFormats()

$ws.Range("E20").Value = "  +0.81%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "1.966.68"
$cell.ClearFormats()

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.404"
$cell.ClearFormats()

$ws.Range("E22").Value = "  +1.52%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "8.711"
$cell.ClearFormats()

$ws.Range("E23").Value = "  -0.33%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "5.400"
$cell.ClearFormats()

$ws.Range("E24").Value = "  +3.15%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "136.70"
$cell.ClearFormats()

$ws.Range("E25").Value = "  -2.13%  "

$ws.Range("E26").Value = "  -1.25%  "

$ws.Range("E27").Value = "  -0.33%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.776"
$cell.ClearFormats()

$ws.Range("E28").Value = "  -1.73%  "

$ws.Range("E29").Value = "  +0.50%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "3.867"
$cell.ClearFormats()

$ws.Range("E30").Value = "  +2.62%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.08212"
$cell.ClearFormats()

$ws.Range("E31").Value = "  -1.13%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.633"
$cell.ClearFormats()

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.04680"
$cell.ClearFormats()

$ws.Range("E34").Value = "  +0.49%  "

$ws.Range("E35").Value = "  -1.48%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.6195"
$cell.ClearFormats()

$ws.Range("E36").Value = "  -2.00%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.750"
$cell.ClearFormats()

$ws.Range("E37").Value = "  +1.58%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01608"
$cell.ClearFormats()

$ws.Range("E38").Value = "  -1.05%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.925"
$cell.ClearFormats()

$ws.Range("E39").Value = "  -2.84%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.9997"
$cell.ClearFormats()

$ws.Range("E40").Value = "  +0.03%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "100.46"
$cell.ClearFormats()

$ws.Range("E41").Value = "  -1.64%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.3925"
$cell.ClearFormats()

$ws.Range("E42").Value = "  -0.21%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.7593"
$cell.ClearFormats()

$ws.Range("E43").Value = "  +0.72%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "5.008"
$cell.ClearFormats()

$ws.Range("E44").Value = "  -1.34%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.1148"
$cell.ClearFormats()

$ws.Range("E45").Value = "  -0.45%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "6.297"
$cell.ClearFormats()

$ws.Range("E46").Value = "  -0.88%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "55.58"
$cell.ClearFormats()

$ws.Range("E47").Value = "  +1.64%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.05272"
$cell.ClearFormats()

$ws.Range("E48").Value = "  -1.33%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "30.65"
$cell.ClearFormats()

$ws.Range("E49").Value = "  -1.32%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "7.535"
$cell.ClearFormats()

$ws.Range("E50").Value = "  -1.15%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.3434"
$cell.ClearFormats()

$ws.Range("E51").Value = "  -1.33%  "
